$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 7299
    $ws.Range("F4").Value = 5567
    $ws.Range("F5").Value = 83
    $ws.Range("F6").Value = 173
    $ws.Range("F9").Value = 107
    $ws.Range("F11").Value = 108
    $ws.Range("F12").Value = 203

    $ws.Range("C13").Value = "【大会员提前抢】合肥·第十三届合肥次元之门动漫游戏博览会-多多poi&Mace专场"
    $ws.Range("F13").Value = 49

    $ws.Range("F14").Value = 649
    $ws.Range("F15").Value = 337
    $ws.Range("F19").Value = 44
}
